# Update BIBI monthly revenue (faturamento_mensal) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (month 1)
$ws.Range("C2").Value = 68977.03999999999
$ws.Range("D2").Value = 243459.36
$ws.Range("E2").Value = 353395.61
$ws.Range("F2").Value = 362289.33
$ws.Range("G2").Value = 416610.57
$ws.Range("H2").Value = 369665.08
$ws.Range("I2").Value = 587582.36

# Row 3 (month 2)
$ws.Range("C3").Value = 54846.21
$ws.Range("D3").Value = 203458.99
$ws.Range("E3").Value = 250171.15
$ws.Range("F3").Value = 338247.04
$ws.Range("G3").Value = 368853.4
$ws.Range("H3").Value = 322821.24
$ws.Range("I3").Value = 763210.77

# Row 4 (month 3)
$ws.Range("C4").Value = 54958
$ws.Range("D4").Value = 185876.61
$ws.Range("E4").Value = 346618.91
$ws.Range("F4").Value = 374906.45
$ws.Range("G4").Value = 402375.47
$ws.Range("H4").Value = 394663.49
$ws.Range("I4").Value = 670620.61

# Row 5 (month 4)
$ws.Range("C5").Value = 62793
$ws.Range("D5").Value = 176355.8
$ws.Range("E5").Value = 314021.11
$ws.Range("F5").Value = 381633.47
$ws.Range("G5").Value = 336910.07
$ws.Range("H5").Value = 450719.69
$ws.Range("I5").Value = 511614.46

# Row 6 (month 5)
$ws.Range("C6").Value = 56618.7
$ws.Range("D6").Value = 227009.89
$ws.Range("E6").Value = 361191.04
$ws.Range("F6").Value = 398659.35
$ws.Range("G6").Value = 371532.61
$ws.Range("H6").Value = 513540.06

# Row 7 (month 6)
$ws.Range("C7").Value = 73560.60000000001
$ws.Range("D7").Value = 231451.45
$ws.Range("E7").Value = 405516
$ws.Range("F7").Value = 372513.52
$ws.Range("G7").Value = 345008.16
$ws.Range("H7").Value = 510962.61

# Row 8 (month 7)
$ws.Range("B8").Value = 74343.75
$ws.Range("C8").Value = 98652.3
$ws.Range("D8").Value = 262869.11
$ws.Range("E8").Value = 543603.28
$ws.Range("F8").Value = 362600.57
$ws.Range("G8").Value = 391131.85
$ws.Range("H8").Value = 647783.67

# Row 9 (month 8)
$ws.Range("B9").Value = 106882
$ws.Range("C9").Value = 178007.6
$ws.Range("D9").Value = 299331.4
$ws.Range("E9").Value = 449321.75
$ws.Range("F9").Value = 360070.83
$ws.Range("G9").Value = 324052.02
$ws.Range("H9").Value = 641898.98

# Row 10 (month 9)
$ws.Range("B10").Value = 82944.5
$ws.Range("C10").Value = 158716.89
$ws.Range("D10").Value = 253588.89
$ws.Range("E10").Value = 458913.37
$ws.Range("F10").Value = 507156.61
$ws.Range("G10").Value = 331199.82
$ws.Range("H10").Value = 535002.39

# Row 11 (month 10)
$ws.Range("B11").Value = 79805.75
$ws.Range("C11").Value = 274966.43
$ws.Range("D11").Value = 219105.59
$ws.Range("E11").Value = 418903.43
$ws.Range("F11").Value = 638193.91
$ws.Range("G11").Value = 361330.03
$ws.Range("H11").Value = 656485.46

# Row 12 (month 11)
$ws.Range("B12").Value = 92135.5
$ws.Range("C12").Value = 225669.34
$ws.Range("D12").Value = 225582.21
$ws.Range("E12").Value = 518320.37
$ws.Range("F12").Value = 487930.04
$ws.Range("G12").Value = 333870.96
$ws.Range("H12").Value = 628940.73

# Row 13 (month 12)
$ws.Range("B13").Value = 125828
$ws.Range("C13").Value = 278306.18
$ws.Range("D13").Value = 260876.34
$ws.Range("E13").Value = 430138.51
$ws.Range("F13").Value = 540802.79
$ws.Range("G13").Value = 487200.04
$ws.Range("H13").Value = 830722.87
